# Change the "name" (col I) and "model_type" (col J) values from "bridge"
# to "link" for a specific set of rows in Sheet1, as described in the
# commit "132 bridges to link".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(3,7,35,103,106,109,113,118,126,130,134,144,147,152,153,157,183,188,195,202,208,215,223,232,235,254,314,315,321,327,334,335,340,349,359,367,375,394,396,399,400,401,422,424,430,431,455,456,505,529,539,540,544,598,599,600,609,611,648,649,651,653,654,660,662,664,690,692,709,755,788,802,804,806,816,842,844,871,872,874,875,887,888,890,892,903,904,1025,1026,1027,1028,1035,1037,1058,1059,1063,1064,1066,1285,1287,1289,1291,1292,1294,1295,1296,1298,1300,1301,1303,1304,1305,1306,1308,1309,1312,1313,1314,1316,1317,1318,1319,1321,1322,1324,1327,1330,1331,1333,1336,1337)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "link"
    $ws.Cells.Item($r, 10).Value = "link"
}
